$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D column: numeric-looking strings that Excel would otherwise reinterpret ---
# as real numbers (e.g. "9.20" -> 9.2, "1.00" -> 1). Force each cell to Text
# format before assigning, then reset the format back to General so the cell
# ends up with no special style attached (same as the rest of the sheet).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D5").Value = "322.68"
$ws.Range("D6").Value = "105.12"
$ws.Range("D7").Value = "0.527"
$ws.Range("D10").Value = "38.15"
$ws.Range("D13").Value = "18.33"
$ws.Range("D19").Value = "12.80"
$ws.Range("D22").Value = "70.74"
$ws.Range("D23").Value = "251.55"
$ws.Range("D25").Value = "2.57"
$ws.Range("D26").Value = "26.19"
$ws.Range("D29").Value = "2.21"
$ws.Range("D30").Value = "35.17"
$ws.Range("D32").Value = "49.41"
$ws.Range("D33").Value = "19.80"
$ws.Range("D34").Value = "5.36"
$ws.Range("D36").Value = "1.00"
$ws.Range("D38").Value = "4.64"
$ws.Range("D42").Value = "122.16"
$ws.Range("D43").Value = "21.37"
$ws.Range("D48").Value = "1.81"
$ws.Range("D49").Value = "9.20"
$ws.Range("D51").Value = "79.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"

# --- D column: values that are already unambiguously text (multiple dots) ---
$ws.Range("D2").Value = "47.582.22"
$ws.Range("D3").Value = "2.492.77"
$ws.Range("D15").Value = "2.883.03"
$ws.Range("D16").Value = "2.491.68"
$ws.Range("D18").Value = "47.459.40"
$ws.Range("D45").Value = "1.967.65"

# --- B/C columns: row 48/49 content swap (FraxShare <-> Stacks) ---
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

# --- E column: Volume(1h) percentage strings (padded with spaces, stay text) ---
$ws.Range("E2").Value = "  +4.72%  "
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  +7.08%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("E18").Value = "  +4.66%  "
$ws.Range("E19").Value = "  +4.51%  "
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("E22").Value = "  +2.64%  "
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("E24").Value = "  +5.80%  "
$ws.Range("E25").Value = "  +3.19%  "
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +4.63%  "
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("E30").Value = "  +6.83%  "
$ws.Range("E31").Value = "  +7.40%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  -2.77%  "
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  +5.66%  "
$ws.Range("E39").Value = "  +4.54%  "
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  +2.07%  "
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("E43").Value = "  +4.03%  "
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("E50").Value = "  +11.50%  "
$ws.Range("E51").Value = "  +3.81%  "
